# Mise a jour des activite et du cours Excel
#
# Renames the two worksheets ("Activite 13" -> "Activite 17",
# "Activite 14" -> "Activite 18") and refreshes the header/footer
# font-style token used in the page setup ("Regular" -> "Normal") on
# both sheets.

$wb = $excel.ActiveWorkbook

$renames = @{
    "Activité 13" = "Activité 17"
    "Activité 14" = "Activité 18"
}

foreach ($ws in $wb.Worksheets) {
    if ($renames.ContainsKey($ws.Name)) {
        $ws.Name = $renames[$ws.Name]
    }
}

foreach ($ws in $wb.Worksheets) {
    $ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
    $ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
}
